# Add a "Save" column (H) to the sheet, mirroring the style of the
# existing header row (style of G1) for H1, and plain numeric cells
# (style of G2:G18) for H2:H18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell - copy the style used by the other header cells (e.g. G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data values for the new "Save" column, row by row
$saveValues = @{
    2  = 0
    3  = 1
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 1
    15 = 1
    16 = 1
    17 = 0
    18 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
